$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = "https://youtu.be/z1VFQ0a8zMY"
$ws.Range("A6").Value = "Lecture 3&4"

$ws.Range("D8").Select()
